$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.5292044909545227
$ws.Range("C2").Value = 0.678826866322334
$ws.Range("D2").Value = 0.739899615432439
$ws.Range("E2").Value = 0.7653563574070629
$ws.Range("B3").Value = 0.5348675326294656
$ws.Range("C3").Value = 0.6839437062193391
$ws.Range("D3").Value = 0.7431061246380334
$ws.Range("E3").Value = 0.767545930173578
$ws.Range("B4").Value = 0.4753858453573038
$ws.Range("C4").Value = 0.6449467475171032
$ws.Range("D4").Value = 0.7215458874735288
$ws.Range("E4").Value = 0.7489833093671484
$ws.Range("B5").Value = 0.5339142996764678
$ws.Range("C5").Value = 0.6864429141796565
$ws.Range("D5").Value = 0.7436580617940645
$ws.Range("E5").Value = 0.7722547943492591
$ws.Range("B6").Value = 0.5062968331428027
$ws.Range("C6").Value = 0.6673089282778126
$ws.Range("D6").Value = 0.7357746999847482
$ws.Range("E6").Value = 0.766576297593929
$ws.Range("B7").Value = 0.5625672726102521
$ws.Range("C7").Value = 0.7051670157368789
$ws.Range("D7").Value = 0.7585522726615487
$ws.Range("E7").Value = 0.7730594823287936
$ws.Range("B8").Value = 0.2571093819663559
$ws.Range("C8").Value = 0.466464064824418
$ws.Range("D8").Value = 0.6022878936784698
$ws.Range("E8").Value = 0.7053807860549471
$ws.Range("B9").Value = 0.4971649692636496
$ws.Range("C9").Value = 0.6519148713952458
$ws.Range("D9").Value = 0.7261943490557077
$ws.Range("E9").Value = 0.7539730207653869
$ws.Range("B10").Value = 0.5948536767995212
$ws.Range("C10").Value = 0.7341654369609846
$ws.Range("D10").Value = 0.7743276069489899
$ws.Range("E10").Value = 0.790092921855691
$ws.Range("B11").Value = 0.5887919840217942
$ws.Range("C11").Value = 0.7300701472138008
$ws.Range("D11").Value = 0.7720118857993419
$ws.Range("E11").Value = 0.7882143195690925
$ws.Range("B12").Value = 0.5416740563608959
$ws.Range("C12").Value = 0.6980149801332342
$ws.Range("D12").Value = 0.7574339153657651
$ws.Range("E12").Value = 0.7731175937837654
$ws.Range("B13").Value = 0.5790612702473026
$ws.Range("C13").Value = 0.7235535682577257
$ws.Range("D13").Value = 0.7686031626159526
$ws.Range("E13").Value = 0.7854760679862451
